$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Duplicate the "Interests Import Queries" sheet so we end up with two
#    sheets: one for Interests-Events relationships (existing data) and one
#    for the new Interests-Users relationships.
# ---------------------------------------------------------------------------
$srcSheet = $wb.Worksheets.Item("Interests Import Queries")
$srcSheet.Copy($null, $srcSheet)

$eventsSheet = $wb.Worksheets.Item("Interests Import Queries")
$usersSheet  = $wb.Worksheets.Item("Interests Import Queries (2)")

$eventsSheet.Name = "Interests-Events Import Queries"
$usersSheet.Name  = "Interests-Users Import Queries"

# ---------------------------------------------------------------------------
# 2. Populate the new "Interests-Users Import Queries" sheet with the
#    user -> interest relationship rows.
# ---------------------------------------------------------------------------
$E = @(1,1,2,2,2,3,3,10,14,11,11,12,12,17,13,5,6,7,8,9,10,11,12,13,14)
$F = @(41,46,41,47,45,48,42,45,41,41,42,41,46,41,41,42,49,49,50,50,50,50,50,50,50)

for ($i = 0; $i -lt $E.Length; $i++) {
    $row = $i + 2
    $usersSheet.Range("E$row").Value = $E[$i]
    $usersSheet.Range("F$row").Value = $F[$i]
    $usersSheet.Range("G$row").Value = "INTERESTED_IN"
}

# Rebuild the H-column relationship-query formulas as two shared groups,
# matching rows 2:21 and the newly appended rows 22:26.
$usersSheet.Range("H2:H21").Formula = '=((((("start n1=node("&E2)&"),n2=node(")&F2)&") create n1-[:")&G2)&"]->n2;"'
$usersSheet.Range("H22:H26").Formula = '=((((("start n1=node("&E22)&"),n2=node(")&F22)&") create n1-[:")&G22)&"]->n2;"'

# The relationship-type column is now wider to fit "INTERESTED_IN".
$usersSheet.Columns.Item(7).ColumnWidth = 15

# ---------------------------------------------------------------------------
# 3. Cosmetic selection / view bookkeeping to mirror the authored workbook.
# ---------------------------------------------------------------------------
$eventsSheet.Activate()
$eventsSheet.Range("F21").Select()

$usersSheet.Activate()
$usersSheet.Range("B2").Select()

# "Events Import Queries" ends up as the active tab / selection in the
# saved workbook.
$eventsImportSheet = $wb.Worksheets.Item("Events Import Queries")
$eventsImportSheet.Activate()
$eventsImportSheet.Range("C23").Select()
